$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Percent" header (previously in E1) is removed; "Unit" (previously in F1)
# takes its place in E1, and F1 becomes blank.
$ws.Range("E1").Value = "Unit"
$ws.Range("F1").Value = $null

# Update the active selection on the sheet to match the saved view.
$ws.Range("D5").Select()

# Best-effort: also reposition/resize the saved workbook window to match the
# author's view (no-op in hosts that don't persist window geometry).
$excel.ActiveWindow.Left = 1140
$excel.ActiveWindow.Top = 1140
$excel.ActiveWindow.Width = 14400
$excel.ActiveWindow.Height = 7270
